$d = $word.ActiveDocument

# Remove the two trailing informational paragraphs ("Ver no Jupiter ...",
# "(c) 2020 ...") plus the blank paragraph that separated them from the
# "LOM3212: ..." requisito line, leaving that paragraph directly followed
# by the page-break paragraph that used to come after them.

for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "*Ver no Jupiter*") {
        $p.Range.Delete() | Out-Null
    }
}

for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "*Contact: luizeleno@usp.br*") {
        $p.Range.Delete() | Out-Null
    }
}

# Now remove the now-orphaned blank paragraph that sat right after the
# "LOM3212: Fenômenos de Transporte A (Requisito)" paragraph.
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "*LOM3212: Fenômenos de Transporte A (Requisito)*") {
        $next = $d.Paragraphs.Item($i + 1)
        if ($next.Range.Text.Trim().Length -eq 0) {
            $next.Range.Delete() | Out-Null
        }
        break
    }
}
